# Doing Updates for Financials
# Apply the restated financial figures (yearly financials refresh) to the YGEHY sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YGEHY")

$ws.Range("D8").Value = 1241300
$ws.Range("E8").Value = 1243100
$ws.Range("F8").Value = 1479000
$ws.Range("G8").Value = 1918600
$ws.Range("H8").Value = 1991400
$ws.Range("I8").Value = 1690700
$ws.Range("J8").Value = 2178400
$ws.Range("D9").Value = 1217800
$ws.Range("E9").Value = 1087100
$ws.Range("F9").Value = 1391800
$ws.Range("G9").Value = 3172700
$ws.Range("H9").Value = 1846100
$ws.Range("I9").Value = 1745400
$ws.Range("J9").Value = 1941200
$ws.Range("D10").Value = 23400
$ws.Range("E10").Value = 156000
$ws.Range("F10").Value = 87200
$ws.Range("G10").Value = -1254200
$ws.Range("H10").Value = 145200
$ws.Range("I10").Value = -54700
$ws.Range("J10").Value = 237200
$ws.Range("D12").Value = 19600
$ws.Range("E12").Value = 20100
$ws.Range("F12").Value = 55900
$ws.Range("G12").Value = 80700
$ws.Range("H12").Value = 39200
$ws.Range("I12").Value = 23400
$ws.Range("J12").Value = 37200
$ws.Range("D14").Value = 303100
$ws.Range("E14").Value = 189600
$ws.Range("F14").Value = 564600
$ws.Range("I14").Value = 29800
$ws.Range("J14").Value = 370500
$ws.Range("D15").Value = 5300
$ws.Range("E15").Value = 7200
$ws.Range("F15").Value = 13100
$ws.Range("G15").Value = 14300
$ws.Range("H15").Value = 13400
$ws.Range("I15").Value = 19900
$ws.Range("J15").Value = 18100
$ws.Range("D17").Value = 1678000
$ws.Range("E17").Value = 1493100
$ws.Range("F17").Value = 2106500
$ws.Range("G17").Value = 1950500
$ws.Range("H17").Value = 2157400
$ws.Range("I17").Value = 2065200
$ws.Range("J17").Value = 2570600
$ws.Range("D18").Value = -436800
$ws.Range("E18").Value = -250100
$ws.Range("F18").Value = -627500
$ws.Range("G18").Value = -31900
$ws.Range("H18").Value = -166000
$ws.Range("I18").Value = -374500
$ws.Range("J18").Value = -392200
$ws.Range("D20").Value = 23200
$ws.Range("E20").Value = 38500
$ws.Range("F20").Value = 5700
$ws.Range("G20").Value = -12300
$ws.Range("H20").Value = 9800
$ws.Range("I20").Value = 3400
$ws.Range("J20").Value = -19200
$ws.Range("D21").Value = -339800
$ws.Range("E21").Value = -78200
$ws.Range("F21").Value = -445600
$ws.Range("G21").Value = 165400
$ws.Range("H21").Value = 36800
$ws.Range("I21").Value = -193300
$ws.Range("J21").Value = -268800
$ws.Range("D22").Value = 98500
$ws.Range("E22").Value = 98400
$ws.Range("F22").Value = 145000
$ws.Range("G22").Value = 150800
$ws.Range("H22").Value = 144200
$ws.Range("I22").Value = 133100
$ws.Range("J22").Value = 93000
$ws.Range("D23").Value = -512000
$ws.Range("E23").Value = -310000
$ws.Range("F23").Value = -766800
$ws.Range("G23").Value = -195000
$ws.Range("H23").Value = -300400
$ws.Range("I23").Value = -504200
$ws.Range("J23").Value = -504500
$ws.Range("F24").Value = 108500
$ws.Range("G24").Value = 13300
$ws.Range("H24").Value = 4600
$ws.Range("I24").Value = -30500
$ws.Range("J24").Value = -19800
$ws.Range("D26").Value = -512100
$ws.Range("E26").Value = -311900
$ws.Range("F26").Value = -875300
$ws.Range("G26").Value = -208300
$ws.Range("H26").Value = -305000
$ws.Range("I26").Value = -473700
$ws.Range("J26").Value = -484700
$ws.Range("D27").Value = -492400
$ws.Range("E27").Value = -311300
$ws.Range("F27").Value = -831200
$ws.Range("G27").Value = -192900
$ws.Range("H27").Value = -288600
$ws.Range("I27").Value = -454800
$ws.Range("J27").Value = -476200
$ws.Range("D32").Value = -23200
$ws.Range("E32").Value = -38500
$ws.Range("F32").Value = -5700
$ws.Range("G32").Value = 12300
$ws.Range("H32").Value = -9800
$ws.Range("I32").Value = -3400
$ws.Range("J32").Value = 19200
$ws.Range("D33").Value = -492400
$ws.Range("E33").Value = -311300
$ws.Range("F33").Value = -831200
$ws.Range("G33").Value = -192900
$ws.Range("H33").Value = -288600
$ws.Range("I33").Value = -454800
$ws.Range("J33").Value = -476200
$ws.Range("D35").Value = -492400
$ws.Range("E35").Value = -311300
$ws.Range("F35").Value = -831200
$ws.Range("G35").Value = -192900
$ws.Range("H35").Value = -288600
$ws.Range("I35").Value = -454800
$ws.Range("J35").Value = -476200
$ws.Range("D41").Value = 56100
$ws.Range("E41").Value = 75200
$ws.Range("F41").Value = 235600
$ws.Range("G41").Value = 158700
$ws.Range("H41").Value = 164100
$ws.Range("I41").Value = 452900
$ws.Range("J41").Value = 620500
$ws.Range("D42").Value = 5500
$ws.Range("D43").Value = 561600
$ws.Range("E43").Value = 557600
$ws.Range("F43").Value = 1123100
$ws.Range("G43").Value = 770200
$ws.Range("H43").Value = 807700
$ws.Range("I43").Value = 1279900
$ws.Range("J43").Value = 516100
$ws.Range("D44").Value = 168200
$ws.Range("E44").Value = 195100
$ws.Range("F44").Value = 440600
$ws.Range("G44").Value = 311500
$ws.Range("H44").Value = 321300
$ws.Range("I44").Value = 749000
$ws.Range("J44").Value = 396400
$ws.Range("D45").Value = 153800
$ws.Range("E45").Value = 156000
$ws.Range("F45").Value = 447500
$ws.Range("G45").Value = 440300
$ws.Range("H45").Value = 392300
$ws.Range("I45").Value = 411000
$ws.Range("J45").Value = 357500
$ws.Range("D46").Value = 945200
$ws.Range("E46").Value = 983900
$ws.Range("F46").Value = 1247100
$ws.Range("G46").Value = 1680600
$ws.Range("H46").Value = 1685400
$ws.Range("I46").Value = 1649300
$ws.Range("J46").Value = 1890500
$ws.Range("D47").Value = 38100
$ws.Range("E47").Value = 62700
$ws.Range("F47").Value = 79900
$ws.Range("G47").Value = 62800
$ws.Range("H47").Value = 32100
$ws.Range("D48").Value = 398900
$ws.Range("E48").Value = 724100
$ws.Range("F48").Value = 2032200
$ws.Range("G48").Value = 1797400
$ws.Range("H48").Value = 1948400
$ws.Range("I48").Value = 3923400
$ws.Range("J48").Value = 1838700
$ws.Range("D49").Value = 67000
$ws.Range("E49").Value = 68400
$ws.Range("F49").Value = 78400
$ws.Range("G49").Value = 98300
$ws.Range("H49").Value = 113100
$ws.Range("I49").Value = 124000
$ws.Range("J49").Value = 94000
$ws.Range("D52").Value = 85700
$ws.Range("E52").Value = 164400
$ws.Range("F52").Value = 222000
$ws.Range("G52").Value = 384100
$ws.Range("H52").Value = 432900
$ws.Range("I52").Value = 304200
$ws.Range("J52").Value = 251900
$ws.Range("D54").Value = 1535000
$ws.Range("E54").Value = 2003500
$ws.Range("F54").Value = 2618000
$ws.Range("G54").Value = 4023200
$ws.Range("H54").Value = 4073800
$ws.Range("I54").Value = 4029900
$ws.Range("J54").Value = 4078800
$ws.Range("D57").Value = 391100
$ws.Range("E57").Value = 425200
$ws.Range("F57").Value = 667400
$ws.Range("G57").Value = 922000
$ws.Range("H57").Value = 904900
$ws.Range("I57").Value = 598700
$ws.Range("J57").Value = 477900
$ws.Range("D58").Value = 1560800
$ws.Range("E58").Value = 2682500
$ws.Range("F58").Value = 1365200
$ws.Range("G58").Value = 1511400
$ws.Range("H58").Value = 2000700
$ws.Range("I58").Value = 1116900
$ws.Range("J58").Value = 1241400
$ws.Range("D59").Value = 428000
$ws.Range("E59").Value = 322300
$ws.Range("F59").Value = 578600
$ws.Range("G59").Value = 249400
$ws.Range("H59").Value = 213700
$ws.Range("I59").Value = 378200
$ws.Range("J59").Value = 202600
$ws.Range("D60").Value = 2379900
$ws.Range("E60").Value = 2092100
$ws.Range("F60").Value = 2324200
$ws.Range("G60").Value = 2682800
$ws.Range("H60").Value = 2122600
$ws.Range("I60").Value = 1920600
$ws.Range("J60").Value = 1921900
$ws.Range("D61").Value = 170600
$ws.Range("E61").Value = 415800
$ws.Range("F61").Value = 401600
$ws.Range("G61").Value = 722500
$ws.Range("H61").Value = 1232500
$ws.Range("I61").Value = 1185600
$ws.Range("J61").Value = 869400
$ws.Range("D62").Value = 541700
$ws.Range("E62").Value = 564300
$ws.Range("F62").Value = 1033300
$ws.Range("G62").Value = 424700
$ws.Range("H62").Value = 402900
$ws.Range("I62").Value = 424300
$ws.Range("J62").Value = 240700
$ws.Range("D66").Value = 3239100
$ws.Range("E66").Value = 3242500
$ws.Range("F66").Value = 3499500
$ws.Range("G66").Value = 4055400
$ws.Range("H66").Value = 3998400
$ws.Range("I66").Value = 3678200
$ws.Range("J66").Value = 3312800
$ws.Range("D72").Value = -2782200
$ws.Range("E72").Value = -2289700
$ws.Range("F72").Value = -1966900
$ws.Range("G72").Value = -1135400
$ws.Range("H72").Value = -942500
$ws.Range("I72").Value = -654000
$ws.Range("J72").Value = -199200
$ws.Range("D76").Value = -1704100
$ws.Range("E76").Value = -1239000
$ws.Range("F76").Value = -881500
$ws.Range("G76").Value = -32300
$ws.Range("H76").Value = 75400
$ws.Range("I76").Value = 351700
$ws.Range("J76").Value = 766000
$ws.Range("D81").Value = -492400
$ws.Range("E81").Value = -311300
$ws.Range("F81").Value = -831200
$ws.Range("G81").Value = -192900
$ws.Range("H81").Value = -288600
$ws.Range("I81").Value = -454800
$ws.Range("J81").Value = -476200
$ws.Range("D83").Value = 73600
$ws.Range("E83").Value = 133100
$ws.Range("F83").Value = 175900
$ws.Range("G83").Value = 209300
$ws.Range("H83").Value = 192700
$ws.Range("I83").Value = 177500
$ws.Range("J83").Value = 142400
$ws.Range("D89").Value = 15000
$ws.Range("E89").Value = -61700
$ws.Range("F89").Value = 146100
$ws.Range("G89").Value = 64400
$ws.Range("H89").Value = 52300
$ws.Range("I89").Value = -305100
$ws.Range("J89").Value = 50600
$ws.Range("D91").Value = -46200
$ws.Range("E91").Value = -36500
$ws.Range("F91").Value = -50100
$ws.Range("G91").Value = -45200
$ws.Range("H91").Value = -176800
$ws.Range("I91").Value = -289800
$ws.Range("J91").Value = -718100
$ws.Range("D94").Value = -34600
$ws.Range("E94").Value = -59100
$ws.Range("F94").Value = 209800
$ws.Range("G94").Value = -46900
$ws.Range("H94").Value = -99500
$ws.Range("I94").Value = -300900
$ws.Range("J94").Value = -816900
$ws.Range("D100").Value = 4400
$ws.Range("E100").Value = 4300
$ws.Range("F100").Value = -333900
$ws.Range("G100").Value = -27100
$ws.Range("H100").Value = -93400
$ws.Range("I100").Value = 297400
$ws.Range("J100").Value = 527600
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = 7500
$ws.Range("G101").Value = 4200
$ws.Range("H101").Value = -4700
$ws.Range("I101").Value = -2500
$ws.Range("J101").Value = -9800
$ws.Range("D102").Value = -19100
$ws.Range("E102").Value = -108900
$ws.Range("F102").Value = 25500
$ws.Range("G102").Value = -5400
$ws.Range("H102").Value = -145400
$ws.Range("I102").Value = -311100
$ws.Range("J102").Value = -248600
